# Auto update Excel log
# Appends new sensor-log rows to the "PIR" and "Humidity" sheets.
# Columns: A=Date, B=Timestamp, C=Hour, D=Location, E=Value, F=Status
# Date/time-like and percentage-like text values are forced to Text
# number format *before* assignment so Excel stores them as literal
# strings instead of auto-converting them into date/time/percent numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: append rows 107-120
# ---------------------------------------------------------------------
$pirSheet = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-01-30","18:27:30","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:30","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:34","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:38","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:44","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:49","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:54","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:27:59","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:28:04","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:28:09","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:28:14","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:28:19","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:28:24","18:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","18:28:29","18:00","Bathroom","No Motion","Inactive")
)

$r = 107
foreach ($row in $pirRows) {
    $pirSheet.Cells.Item($r, 1).NumberFormat = "@"
    $pirSheet.Cells.Item($r, 1).Value = $row[0]

    $pirSheet.Cells.Item($r, 2).NumberFormat = "@"
    $pirSheet.Cells.Item($r, 2).Value = $row[1]

    $pirSheet.Cells.Item($r, 3).NumberFormat = "@"
    $pirSheet.Cells.Item($r, 3).Value = $row[2]

    $pirSheet.Cells.Item($r, 4).Value = $row[3]
    $pirSheet.Cells.Item($r, 5).Value = $row[4]
    $pirSheet.Cells.Item($r, 6).Value = $row[5]

    $r++
}

# ---------------------------------------------------------------------
# Humidity sheet: append rows 70-78
# ---------------------------------------------------------------------
$humiditySheet = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("2026-01-30","18:27:30","18:00","Bathroom","86.6%","Active"),
    @("2026-01-30","18:27:31","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:27:44","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:27:49","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:27:54","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:28:04","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:28:09","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:28:24","18:00","Bathroom","86.5%","Active"),
    @("2026-01-30","18:28:29","18:00","Bathroom","86.5%","Active")
)

$r = 70
foreach ($row in $humidityRows) {
    $humiditySheet.Cells.Item($r, 1).NumberFormat = "@"
    $humiditySheet.Cells.Item($r, 1).Value = $row[0]

    $humiditySheet.Cells.Item($r, 2).NumberFormat = "@"
    $humiditySheet.Cells.Item($r, 2).Value = $row[1]

    $humiditySheet.Cells.Item($r, 3).NumberFormat = "@"
    $humiditySheet.Cells.Item($r, 3).Value = $row[2]

    $humiditySheet.Cells.Item($r, 4).Value = $row[3]

    # Value column looks like a percentage ("86.6%"); force text so it
    # is stored as the literal string, not a numeric percent.
    $humiditySheet.Cells.Item($r, 5).NumberFormat = "@"
    $humiditySheet.Cells.Item($r, 5).Value = $row[4]

    $humiditySheet.Cells.Item($r, 6).Value = $row[5]

    $r++
}

"Added $($pirRows.Count) rows to PIR and $($humidityRows.Count) rows to Humidity"
